$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.955.16"
$ws.Range("D3").Value = "2.417.42"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'552.84"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").Value = "'137.83"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.575"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").Value = "'5.75"
$ws.Range("E10").Value = "  +4.66%  "
$ws.Range("D11").Value = "'0.359"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "'24.86"
$ws.Range("E13").Value = "  +4.81%  "
$ws.Range("D14").Value = "2.847.41"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "59.886.99"
$ws.Range("E15").Value = "  +3.70%  "
$ws.Range("D17").Value = "2.420.66"
$ws.Range("E17").Value = "  +3.49%  "
$ws.Range("D18").Value = "'11.33"
$ws.Range("E18").Value = "  +6.20%  "
$ws.Range("D19").Value = "'4.38"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").Value = "'331.26"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "'6.76"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'65.02"
$ws.Range("E23").Value = "  +3.61%  "
$ws.Range("D24").Value = "'0.170"
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("D25").Value = "'8.57"
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "0.0₃0780"
$ws.Range("E28").Value = "  +6.24%  "
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("D31").Value = "'169.24"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'18.71"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("B33").Value = "SuiNetwork"
$ws.Range("C33").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D33").Value = "'1.04"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'1.29"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "'4.20"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "'39.46"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.417"
$ws.Range("E40").Value = "  +10.81%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'318.94"
$ws.Range("E41").Value = "  +10.96%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "'139.47"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").Value = "'0.0960"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D46").Value = "'19.52"
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.573"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").Value = "'0.408"
$ws.Range("E48").Value = "  +7.18%  "
$ws.Range("D49").Value = "'0.0226"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "'17.68"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "'11.05"
$ws.Range("E51").Value = "  -0.23%  "

Write-Host "Updated cryptos list"
